# Auto-generated Excel COM-interop script
# Applies value updates to H..N columns across multiple rows/sheets
$wb = $excel.ActiveWorkbook

# Sheet ALC, Row 28
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H28").Value = 2500
$ws.Range("I28").Value = 0
$ws.Range("J28").Value = 2500
$ws.Range("K28").Value = 0
$ws.Range("L28").Value = 2500
$ws.Range("N28").Value = -3470

# Sheet ALC, Row 40
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H40").Value = 1597.5
$ws.Range("I40").Value = 1516.6666
$ws.Range("J40").Value = 1646
$ws.Range("K40").Value = 1516.6666
$ws.Range("L40").Value = 1646
$ws.Range("M40").Value = -1341.6666

# Sheet ALC, Row 43
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H43").Value = 1330.8
$ws.Range("I43").Value = 899.5
$ws.Range("J43").Value = 1438.625
$ws.Range("K43").Value = 899.5
$ws.Range("L43").Value = 1438.625
$ws.Range("M43").Value = -830.5
$ws.Range("N43").Value = -1576.625

# Sheet ALC, Row 92
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H92").Value = 231.42857
$ws.Range("I92").Value = 231.42857
$ws.Range("J92").Value = 0
$ws.Range("K92").Value = 231.42857
$ws.Range("L92").Value = 0
$ws.Range("M92").Value = 1016.57143

# Sheet ALC, Row 112
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H112").Value = 1064.3043
$ws.Range("I112").Value = 0
$ws.Range("J112").Value = 1064.3043
$ws.Range("K112").Value = 0
$ws.Range("L112").Value = 3192.9129
$ws.Range("N112").Value = -5408.9129

# Sheet ARM, Row 32
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H32").Value = 3975.739
$ws.Range("I32").Value = 2677.946
$ws.Range("J32").Value = 9311.111000000001
$ws.Range("K32").Value = 2677.946
$ws.Range("L32").Value = 9311.111000000001
$ws.Range("M32").Value = -2390.946

# Sheet ARM, Row 61
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H61").Value = 2948.7273
$ws.Range("I61").Value = 2073.2
$ws.Range("J61").Value = 4824.857
$ws.Range("K61").Value = 2073.2
$ws.Range("L61").Value = 4824.857
$ws.Range("M61").Value = -1861.2
$ws.Range("N61").Value = -5248.857

# Sheet ARM, Row 136
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H136").Value = 2948.7273
$ws.Range("I136").Value = 2073.2
$ws.Range("J136").Value = 4824.857
$ws.Range("K136").Value = 6219.599999999999
$ws.Range("L136").Value = 14474.571
$ws.Range("M136").Value = -3669.599999999999
$ws.Range("N136").Value = -19574.571

# Sheet BSM, Row 86
$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H86").Value = 667499.7
$ws.Range("I86").Value = 1500
$ws.Range("J86").Value = 1000499.5
$ws.Range("K86").Value = 1500
$ws.Range("L86").Value = 1000499.5
$ws.Range("M86").Value = -377
$ws.Range("N86").Value = -1002745.5

# Sheet BSM, Row 89
$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H89").Value = 667499.7
$ws.Range("I89").Value = 1500
$ws.Range("J89").Value = 1000499.5
$ws.Range("K89").Value = 7500
$ws.Range("L89").Value = 5002497.5
$ws.Range("M89").Value = -1884
$ws.Range("N89").Value = -5013729.5

# Sheet BSM, Row 97
$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H97").Value = 5681.2
$ws.Range("I97").Value = 5681.2
$ws.Range("J97").Value = 0
$ws.Range("K97").Value = 5681.2
$ws.Range("L97").Value = 0
$ws.Range("M97").Value = -4690.2

# Sheet BSM, Row 98
$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H98").Value = 0
$ws.Range("I98").Value = 0
$ws.Range("J98").Value = 0
$ws.Range("K98").Value = 0
$ws.Range("L98").Value = 0
$ws.Range("N98").ClearContents()

# Sheet BSM, Row 99
$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H99").Value = 1500
$ws.Range("I99").Value = 1500
$ws.Range("J99").Value = 0
$ws.Range("K99").Value = 1500
$ws.Range("L99").Value = 0
$ws.Range("M99").Value = -2
$ws.Range("N99").ClearContents()

# Sheet BSM, Row 105
$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H105").Value = 2476.0952
$ws.Range("I105").Value = 2489.3684
$ws.Range("J105").Value = 2350
$ws.Range("K105").Value = 2489.3684
$ws.Range("L105").Value = 2350
$ws.Range("M105").Value = -742.3683999999998

# Sheet CRP, Row 31
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H31").Value = 1397.5927
$ws.Range("I31").Value = 876.7
$ws.Range("J31").Value = 1704
$ws.Range("K31").Value = 876.7
$ws.Range("L31").Value = 1704
$ws.Range("M31").Value = -581.7
$ws.Range("N31").Value = -2294

# Sheet CRP, Row 34
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H34").Value = 1397.5927
$ws.Range("I34").Value = 876.7
$ws.Range("J34").Value = 1704
$ws.Range("K34").Value = 876.7
$ws.Range("L34").Value = 1704
$ws.Range("M34").Value = -674.7
$ws.Range("N34").Value = -2108

# Sheet CRP, Row 58
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H58").Value = 2175203.5
$ws.Range("I58").Value = 3953972.2
$ws.Range("J58").Value = 1153
$ws.Range("K58").Value = 3953972.2
$ws.Range("L58").Value = 1153
$ws.Range("M58").Value = -3953769.2

# Sheet CRP, Row 110
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H110").Value = 0
$ws.Range("I110").Value = 0
$ws.Range("J110").Value = 0
$ws.Range("K110").Value = 0
$ws.Range("L110").Value = 0
$ws.Range("N110").ClearContents()

# Sheet CRP, Row 136
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H136").Value = 2175203.5
$ws.Range("I136").Value = 3953972.2
$ws.Range("J136").Value = 1153
$ws.Range("K136").Value = 11861916.6
$ws.Range("L136").Value = 3459
$ws.Range("M136").Value = -11859366.6

# Sheet CUL, Row 64
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H64").Value = 3333.3333
$ws.Range("I64").Value = 0
$ws.Range("J64").Value = 3333.3333
$ws.Range("K64").Value = 0
$ws.Range("L64").Value = 9999.999899999999
$ws.Range("M64").ClearContents()
$ws.Range("N64").Value = -10539.9999

# Sheet CUL, Row 67
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H67").Value = 3333.3333
$ws.Range("I67").Value = 0
$ws.Range("J67").Value = 3333.3333
$ws.Range("K67").Value = 0
$ws.Range("L67").Value = 9999.999899999999
$ws.Range("M67").ClearContents()
$ws.Range("N67").Value = -11871.9999

# Sheet CUL, Row 69
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H69").Value = 2900
$ws.Range("I69").Value = 2750
$ws.Range("J69").Value = 3000
$ws.Range("K69").Value = 8250
$ws.Range("L69").Value = 9000
$ws.Range("M69").Value = -7439
$ws.Range("N69").Value = -10622

# Sheet CUL, Row 72
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H72").Value = 2900
$ws.Range("I72").Value = 2750
$ws.Range("J72").Value = 3000
$ws.Range("K72").Value = 24750
$ws.Range("L72").Value = 27000
$ws.Range("M72").Value = -20694
$ws.Range("N72").Value = -35112

# Sheet CUL, Row 107
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H107").Value = 1546.6154
$ws.Range("I107").Value = 1103.9166
$ws.Range("J107").Value = 1926.0714
$ws.Range("K107").Value = 3311.7498
$ws.Range("L107").Value = 5778.2142
$ws.Range("M107").Value = -1391.7498
$ws.Range("N107").Value = -9618.2142

# Sheet CUL, Row 131
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H131").Value = 7826125.5
$ws.Range("I131").Value = 100000350
$ws.Range("J131").Value = 14750.322
$ws.Range("K131").Value = 300001050
$ws.Range("L131").Value = 44250.966
$ws.Range("M131").Value = -299996010
$ws.Range("N131").Value = -54330.966

# Sheet CUL, Row 134
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H134").Value = 4632
$ws.Range("I134").Value = 4143
$ws.Range("J134").Value = 4998.75
$ws.Range("K134").Value = 12429
$ws.Range("L134").Value = 14996.25
$ws.Range("M134").Value = -7359
$ws.Range("N134").Value = -25136.25

# Sheet GSM, Row 126
$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H126").Value = 1716678
$ws.Range("I126").Value = 2060595
$ws.Range("J126").Value = 169052
$ws.Range("K126").Value = 6181785
$ws.Range("L126").Value = 507156
$ws.Range("M126").Value = -6179315
$ws.Range("N126").Value = -512096

# Sheet LTW, Row 46
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H46").Value = 1100
$ws.Range("I46").Value = 450
$ws.Range("J46").Value = 1316.6666
$ws.Range("K46").Value = 450
$ws.Range("L46").Value = 1316.6666
$ws.Range("M46").Value = -262
$ws.Range("N46").Value = -1692.6666

# Sheet LTW, Row 136
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H136").Value = 3209.1614
$ws.Range("I136").Value = 1855.6364
$ws.Range("J136").Value = 6517.778
$ws.Range("K136").Value = 5566.9092
$ws.Range("L136").Value = 19553.334
$ws.Range("M136").Value = -3016.9092

# Sheet WVR, Row 132
$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H132").Value = 1529.7446
$ws.Range("I132").Value = 1161.8235
$ws.Range("J132").Value = 2492
$ws.Range("K132").Value = 3485.4705
$ws.Range("L132").Value = 7476
$ws.Range("M132").Value = -955.4704999999999
$ws.Range("N132").Value = -12536

# Sheet WVR, Row 136
$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H136").Value = 11114222
$ws.Range("I136").Value = 18521410
$ws.Range("J136").Value = 3439.25
$ws.Range("K136").Value = 55564230
$ws.Range("L136").Value = 10317.75
$ws.Range("M136").Value = -55561680
$ws.Range("N136").Value = -15417.75

# Sheet WVR, Row 137
$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H137").Value = 50000
$ws.Range("I137").Value = 0
$ws.Range("J137").Value = 50000
$ws.Range("K137").Value = 0
$ws.Range("L137").Value = 50000
$ws.Range("N137").Value = -60200

# Sheet WVR, Row 141
$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H141").Value = 70578.75
$ws.Range("I141").Value = 0
$ws.Range("J141").Value = 70578.75
$ws.Range("K141").Value = 0
$ws.Range("L141").Value = 70578.75
